$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-CellText "D2" "28.922.48"
Set-CellText "E2" "  -0.39%  "
Set-CellText "D3" "1.823.14"
Set-CellText "E3" "  -0.51%  "
Set-CellText "D4" "0.9964"
Set-CellText "E4" "  -0.27%  "
Set-CellText "D5" "243.10"
Set-CellText "E5" "  +0.65%  "
Set-CellText "E6" "  -0.04%  "
Set-CellText "D7" "0.9962"
Set-CellText "E7" "  -0.41%  "
Set-CellText "D8" "0.07452"
Set-CellText "E8" "  -1.65%  "
Set-CellText "D9" "0.2928"
Set-CellText "E9" "  +0.28%  "
Set-CellText "E10" "  +0.40%  "
Set-CellText "D11" "0.07668"
Set-CellText "E11" "  +0.29%  "
Set-CellText "D12" "1.824.31"
Set-CellText "E12" "  -0.50%  "
Set-CellText "D13" "4.975"
Set-CellText "E13" "  +0.36%  "
Set-CellText "D14" "0.6649"
Set-CellText "E14" "  -0.10%  "
Set-CellText "D15" "82.83"
Set-CellText "E15" "  +0.45%  "
Set-CellText "D16" "0.000009666"
Set-CellText "E16" "  +2.62%  "
Set-CellText "D17" "6.020"
Set-CellText "E17" "  +0.63%  "
Set-CellText "D18" "28.916.20"
Set-CellText "E18" "  -0.40%  "
Set-CellText "E19" "  +1.53%  "
Set-CellText "D20" "224.67"
Set-CellText "E20" "  -0.30%  "
Set-CellText "D21" "0.9945"
Set-CellText "E21" "  -0.50%  "
Set-CellText "D22" "7.108"
Set-CellText "E22" "  -1.53%  "
Set-CellText "D23" "0.9970"
Set-CellText "E23" "  -0.37%  "
Set-CellText "D24" "159.99"
Set-CellText "E24" "  -0.14%  "
Set-CellText "D25" "0.1404"
Set-CellText "E25" "  +2.76%  "
Set-CellText "D26" "8.473"
Set-CellText "E26" "  +0.59%  "
Set-CellText "D27" "17.85"
Set-CellText "D28" "1.491"
Set-CellText "E28" "  -0.34%  "
Set-CellText "E29" "  +1.07%  "
Set-CellText "D30" "4.046"
Set-CellText "E30" "  +0.28%  "
Set-CellText "D31" "0.05437"
Set-CellText "E31" "  +4.36%  "
Set-CellText "D32" "1.197"
Set-CellText "E32" "  -0.15%  "
Set-CellText "E33" "  -0.05%  "
Set-CellText "D34" "0.7404"
Set-CellText "E34" "  +0.99%  "
Set-CellText "D35" "1.131"
Set-CellText "E35" "  -1.73%  "
Set-CellText "D36" "2.606"
Set-CellText "E36" "  +0.70%  "
Set-CellText "D37" "1.237.86"
Set-CellText "E37" "  -2.76%  "
Set-CellText "D38" "2.737"
Set-CellText "E38" "  -0.86%  "
Set-CellText "D39" "0.01772"
Set-CellText "E39" "  -0.85%  "
Set-CellText "D40" "6.645"
Set-CellText "E40" "  +1.34%  "
Set-CellText "D41" "0.8953"
Set-CellText "E41" "  +0.06%  "
Set-CellText "D42" "0.9946"
Set-CellText "E42" "  -0.58%  "
Set-CellText "D43" "101.17"
Set-CellText "E43" "  -0.53%  "
Set-CellText "D44" "1.972.58"
Set-CellText "E44" "  -0.18%  "
Set-CellText "D45" "64.77"
Set-CellText "E45" "  +0.26%  "
Set-CellText "E46" "  +1.69%  "
Set-CellText "D47" "0.5079"
Set-CellText "E47" "  -0.68%  "
Set-CellText "E48" "  +1.40%  "
Set-CellText "B49" "XinFinNetwork"
Set-CellText "C49" "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
Set-CellText "D49" "0.07345"
Set-CellText "E49" "  +3.54%  "
Set-CellText "B50" "EnergySwap"
Set-CellText "C50" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-CellText "D50" "8.950"
Set-CellText "E50" "  +0.94%  "
Set-CellText "D51" "1.656"
Set-CellText "E51" "  +1.48%  "
